$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 0.9997698664665222
$ws.Range("B1").Value = 1.80894935131073
$ws.Range("C1").Value = 1.916611790657043
$ws.Range("D1").Value = 2.009451150894165
$ws.Range("E1").Value = 1.393824934959412
